$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9180859327316284
$ws.Range("B1").Value = 1.096815347671509
$ws.Range("C1").Value = 1.477943181991577
$ws.Range("D1").Value = 3.742326259613037
$ws.Range("E1").Value = 3.711601972579956
